# Refresh the live crypto snapshot (price + 1h volume %, and occasional rank
# swaps when two coins trade places) pulled in by the scheduled GitHub Action.
# Column D ("Price") is a text column even though most values look numeric -
# a bare numeric-looking Value would be auto-coerced to a Double by Excel, so
# those assignments are prefixed with a literal leading apostrophe, Excel's
# normal "force text" entry method (COM strips the apostrophe itself, leaving
# a clean text value - same as typing it into a cell by hand).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "39.887.26"
$ws.Range("E2").Value = "  +0.16%  "
# Row 3
$ws.Range("D3").Value = "2.211.73"
$ws.Range("E3").Value = "  -0.37%  "
# Row 4
$ws.Range("E4").Value = "  -0.08%  "
# Row 5
$ws.Range("D5").Value = "'291.52"
$ws.Range("E5").Value = "  -0.22%  "
# Row 6
$ws.Range("D6").Value = "'86.71"
$ws.Range("E6").Value = "  +0.56%  "
# Row 7
$ws.Range("D7").Value = "'0.512"
$ws.Range("E7").Value = "  -0.42%  "
# Row 8
$ws.Range("E8").Value = "  -0.09%  "
# Row 9
$ws.Range("D9").Value = "'0.467"
$ws.Range("E9").Value = "  -1.04%  "
# Row 10
$ws.Range("D10").Value = "'30.29"
$ws.Range("E10").Value = "  -1.67%  "
# Row 11
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.0778"
$ws.Range("E11").Value = "  -0.94%  "
# Row 12
$ws.Range("B12").Value = "OKB"
$ws.Range("C12").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D12").Value = "'50.04"
$ws.Range("E12").Value = "  +6.07%  "
# Row 13
$ws.Range("E13").Value = "  +2.68%  "
# Row 14
$ws.Range("D14").Value = "'6.42"
$ws.Range("E14").Value = "  +1.04%  "
# Row 15
$ws.Range("D15").Value = "2.554.62"
$ws.Range("E15").Value = "  -0.27%  "
# Row 16
$ws.Range("D16").Value = "'13.70"
$ws.Range("E16").Value = "  -2.41%  "
# Row 17
$ws.Range("D17").Value = "2.212.25"
$ws.Range("E17").Value = "  -0.35%  "
# Row 18
$ws.Range("D18").Value = "'0.729"
$ws.Range("E18").Value = "  -0.40%  "
# Row 19
$ws.Range("D19").Value = "39.771.08"
$ws.Range("E19").Value = "  -0.05%  "
# Row 20
$ws.Range("D20").Value = "0.0₃0883"
$ws.Range("E20").Value = "  +0.08%  "
# Row 21
$ws.Range("D21").Value = "'11.19"
$ws.Range("E21").Value = "  +1.08%  "
# Row 22
$ws.Range("D22").Value = "'5.73"
$ws.Range("E22").Value = "  -1.31%  "
# Row 23
$ws.Range("D23").Value = "'65.44"
$ws.Range("E23").Value = "  -0.27%  "
# Row 24
$ws.Range("D24").Value = "'236.62"
$ws.Range("E24").Value = "  +0.18%  "
# Row 26
$ws.Range("D26").Value = "'2.45"
$ws.Range("E26").Value = "  -0.67%  "
# Row 27
$ws.Range("D27").Value = "'1.83"
$ws.Range("E27").Value = "  -0.88%  "
# Row 28
$ws.Range("D28").Value = "'23.26"
$ws.Range("E28").Value = "  +2.17%  "
# Row 29
$ws.Range("E29").Value = "  -2.85%  "
# Row 30
$ws.Range("D30").Value = "'9.20"
$ws.Range("E30").Value = "  -0.47%  "
# Row 31
$ws.Range("D31").Value = "'157.11"
$ws.Range("E31").Value = "  +3.62%  "
# Row 32
$ws.Range("D32").Value = "'31.82"
$ws.Range("E32").Value = "  -3.23%  "
# Row 33
$ws.Range("E33").Value = "  +0.02%  "
# Row 34
$ws.Range("D34").Value = "'4.94"
$ws.Range("E34").Value = "  -0.02%  "
# Row 35
$ws.Range("D35").Value = "'0.0709"
$ws.Range("E35").Value = "  -1.20%  "
# Row 36
$ws.Range("D36").Value = "'2.91"
$ws.Range("E36").Value = "  +3.87%  "
# Row 37
$ws.Range("E37").Value = "  -1.90%  "
# Row 38
$ws.Range("E38").Value = "  -0.61%  "
# Row 39
$ws.Range("D39").Value = "'0.0982"
$ws.Range("E39").Value = "  -1.26%  "
# Row 40
$ws.Range("E40").Value = "  +0.29%  "
# Row 41
$ws.Range("D41").Value = "'15.19"
$ws.Range("E41").Value = "  -4.17%  "
# Row 42
$ws.Range("D42").Value = "2.109.36"
$ws.Range("E42").Value = "  +2.06%  "
# Row 43
$ws.Range("D43").Value = "'3.71"
$ws.Range("E43").Value = "  -2.09%  "
# Row 44
$ws.Range("D44").Value = "'0.0268"
$ws.Range("E44").Value = "  +0.14%  "
# Row 45
$ws.Range("E45").Value = "  -0.05%  "
# Row 46
$ws.Range("D46").Value = "'9.87"
$ws.Range("E46").Value = "  -0.68%  "
# Row 47
$ws.Range("D47").Value = "'17.66"
$ws.Range("E47").Value = "  -0.87%  "
# Row 48
$ws.Range("D48").Value = "'2.70"
$ws.Range("E48").Value = "  +3.71%  "
# Row 49
$ws.Range("D49").Value = "2.428.99"
$ws.Range("E49").Value = "  -0.29%  "
# Row 50
$ws.Range("D50").Value = "'1.48"
$ws.Range("E50").Value = "  +2.48%  "
# Row 51
$ws.Range("D51").Value = "'88.28"
$ws.Range("E51").Value = "  -0.85%  "
